# ===========================================================================
# Handback status report regeneration.
#
# A new handed-back file (d332a154-1705-4dfb-b938-0752c3ccafb3.md) showed up
# in this run alongside the previously-tracked file, whose generated GUID
# changed to 4f94797c-8777-4592-8da8-a09997000153 and whose xliff timestamps
# were refreshed. Both files get a row appended to every sheet.
# ===========================================================================

function Set-ExistingHyperlink {
    # Updates the Address/TextToDisplay of a hyperlink already anchored to
    # $cellAddr (e.g. "$A$2") without creating a duplicate link object.
    param($ws, $cellAddr, $newAddress, $newDisplay)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $cellAddr) {
            $h.Address = $newAddress
            $h.TextToDisplay = $newDisplay
            return
        }
    }
}

function Add-NewHyperlink {
    # Adds a brand-new hyperlink for a freshly-populated cell and restyles it
    # to look like the workbook's existing "HyperLink" cells (underlined,
    # Excel's classic link blue).
    param($ws, $cell, $address, $display)
    $ws.Hyperlinks.Add($cell, $address, "", "", $display) | Out-Null
    $cell.Font.Underline = $True
    $cell.Font.Color = 15570276
}

$wb = $excel.ActiveWorkbook

$oldGuid = "35fcf230-f3f3-499f-8195-5edfd46dc5d4"
$newGuid = "4f94797c-8777-4592-8da8-a09997000153"
$addGuid = "d332a154-1705-4dfb-b938-0752c3ccafb3"

$newZhXlfHash = "6f65dfdbcf5662186702c45cb114a6f79077a558"
$addXlfHash   = "eb3e0f523112856901ddbcebe4af66710497475a"

$zhcnRepoUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ed0caea2d5ea9201f2f94aacbae1e4afd78d88ea/e2e"
$dedeRepoUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4e1fc295266d11db85f61884316fb7eff654711b/e2e"
$srcRepoUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6835d3a5a5f723b568be90a71d9959554d4d5777/e2e"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 (existing file) picks up the refreshed GUID + generation date.
$ov.Range("A2").Value = "$newGuid.md"
$ov.Range("B2").Value = "e2e\$newGuid.md"
Set-ExistingHyperlink $ov "`$B`$2" "$srcRepoUrl/$newGuid.md" "e2e\$newGuid.md"
$ov.Range("G2").Value = "2016-08-28 23:00:36"

# Row 3 (new file).
$ov.Range("A3").Value = "$addGuid.md"
$ov.Range("B3").Value = "e2e\$addGuid.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2016-08-28 23:00:36"
Add-NewHyperlink $ov $ov.Range("B3") "$srcRepoUrl/$addGuid.md" "e2e\$addGuid.md"

$ov.Columns.Item(1).ColumnWidth = 39.14

$ov.ListObjects.Item(1).Resize($ov.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 (existing file): GUID + xlf hash + timestamps refreshed.
$zh.Range("A2").Value = "$newGuid.md"
$zh.Range("G2").Value = "$newGuid.$newZhXlfHash.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-28 23:00:32"
$zh.Range("I2").Value = "$newGuid.md"
$zh.Range("J2").Value = "$newGuid.$newZhXlfHash.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-28 23:00:48"
Set-ExistingHyperlink $zh "`$A`$2" "$srcRepoUrl/$newGuid.md" "$newGuid.md"
Set-ExistingHyperlink $zh "`$I`$2" "$zhcnRepoUrl/$newGuid.md" "$newGuid.md"

# Row 3 (new file).
$zh.Range("A3").Value = "$addGuid.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "True"
$zh.Range("G3").Value = "$addGuid.$addXlfHash.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-28 23:00:32"
$zh.Range("I3").Value = "$addGuid.md"
$zh.Range("J3").Value = "$addGuid.$addXlfHash.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-28 23:00:48"
$zh.Range("L3").Value = ""
$zh.Range("M3").Value = "True"
$zh.Range("N3").Value = ""
$zh.Range("O3").Value = "False"
$zh.Range("P3").Value = ""
Add-NewHyperlink $zh $zh.Range("A3") "$srcRepoUrl/$addGuid.md" "$addGuid.md"
Add-NewHyperlink $zh $zh.Range("I3") "$zhcnRepoUrl/$addGuid.md" "$addGuid.md"

$zh.Columns.Item(1).ColumnWidth = 39.14
$zh.Columns.Item(9).ColumnWidth = 39.14

$zh.ListObjects.Item(1).Resize($zh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2 (existing file): GUID + xlf hash + timestamps refreshed.
$de.Range("A2").Value = "$newGuid.md"
$de.Range("G2").Value = "$newGuid.$newZhXlfHash.de-de.xlf"
$de.Range("H2").Value = "2016-08-28 23:00:36"
$de.Range("I2").Value = "$newGuid.md"
$de.Range("J2").Value = "$newGuid.$newZhXlfHash.de-de.xlf"
$de.Range("K2").Value = "2016-08-28 23:00:55"
Set-ExistingHyperlink $de "`$A`$2" "$srcRepoUrl/$newGuid.md" "$newGuid.md"
Set-ExistingHyperlink $de "`$I`$2" "$dedeRepoUrl/$newGuid.md" "$newGuid.md"

# Row 3 (new file).
$de.Range("A3").Value = "$addGuid.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "True"
$de.Range("G3").Value = "$addGuid.$addXlfHash.de-de.xlf"
$de.Range("H3").Value = "2016-08-28 23:00:36"
$de.Range("I3").Value = "$addGuid.md"
$de.Range("J3").Value = "$addGuid.$addXlfHash.de-de.xlf"
$de.Range("K3").Value = "2016-08-28 23:00:55"
$de.Range("L3").Value = ""
$de.Range("M3").Value = "True"
$de.Range("N3").Value = ""
$de.Range("O3").Value = "False"
$de.Range("P3").Value = ""
Add-NewHyperlink $de $de.Range("A3") "$srcRepoUrl/$addGuid.md" "$addGuid.md"
Add-NewHyperlink $de $de.Range("I3") "$dedeRepoUrl/$addGuid.md" "$addGuid.md"

$de.Columns.Item(1).ColumnWidth = 39.14
$de.Columns.Item(9).ColumnWidth = 39.14

$de.ListObjects.Item(1).Resize($de.Range("A1:P3"))
